# Update functions and Data Model (#50)
# Adds a new "Authorship Resource" column (K) to Sheet1, crediting the
# resource's authors on every data row, and refreshes the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$header = "Authorship Resource"
$author = "Daniela Subotic, Noémi Villars-Amberg"

# Header cell
$ws.Range("K1").Value = $header

# Same authorship credit for every data row (rows 2-11)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 11).Value = $author
}

# Resize the new column to fit its content, same as the other bestFit columns
$ws.Columns.Item(11).AutoFit()

# Reflect the new selection left on the sheet after the edit
$ws.Range("K2:K11").Select()
